# Apply a reordering (permutation) of rows 69-82 for columns A, B, D, E, F, G, H, Q, R.
# The row identified by its original "Id" (column A) keeps all of its data, but the
# row gets relocated to a different row number in the sheet. All other columns for
# rows 69-82 (C, I, O, P, S..AY) are identical across the block, so they are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by index) that participate in the permutation.
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18)   # A, B, D, E, F, G, H, Q, R

$firstRow = 69
$lastRow = 82

# Snapshot current ("before") values for each row/column so that writes below
# don't clobber values we still need to read.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping of destination row -> source row (source row's original content
# moves into destination row).
$mapping = @{
    69 = 80
    70 = 75
    71 = 70
    72 = 69
    73 = 74
    74 = 73
    75 = 81
    76 = 79
    77 = 76
    78 = 77
    79 = 71
    80 = 82
    81 = 72
    82 = 78
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
